# future.xlsx edit: fill in results for the games played on Tue, Feb 6 2024
# (Sheet1 rows 136-142), which ripples into Sheet2's rolling-accuracy table
# (row 20 becomes populated, rows 20-37 recalc) and the three charts that
# cache Sheet2!A2:A19 / D2:D19 / F2:F19 / G2:G19 (now A2:A20 etc.).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: complete rows 136-142 (actual scores + derived columns) -------
# columns: row | AwayPts(D) | HomePts(F) | Overtime(G) | Win(I) | Loss(J) | Forecasted(K) | Correct(L) | K-highlight
$games = @(
    @{ Row=136; D=129; F=132; G="NA"; I="Indiana Pacers";          J="Houston Rockets";         K="Houston Rockets";         L="No";  Hi=$false },
    @{ Row=137; D=119; F=107; G="NA"; I="Dallas Mavericks";        J="Brooklyn Nets";            K="Dallas Mavericks";        L="Yes"; Hi=$true  },
    @{ Row=138; D=95;  F=121; G="NA"; I="Miami Heat";              J="Orlando Magic";            K="Orlando Magic";           L="No";  Hi=$false },
    @{ Row=139; D=113; F=123; G="NA"; I="New York Knicks";         J="Memphis Grizzlies";        K="New York Knicks";         L="Yes"; Hi=$true  },
    @{ Row=140; D=123; F=129; G="NA"; I="Chicago Bulls";           J="Minnesota Timberwolves";   K="Minnesota Timberwolves";  L="No";  Hi=$false },
    @{ Row=141; D=117; F=124; G="NA"; I="Utah Jazz";                J="Oklahoma City Thunder";   K="Oklahoma City Thunder";   L="No";  Hi=$false },
    @{ Row=142; D=106; F=114; G="NA"; I="Phoenix Suns";             J="Milwaukee Bucks";         K="Milwaukee Bucks";         L="No";  Hi=$false }
)

foreach ($g in $games) {
    $r = $g.Row

    $ws1.Range("D$r").Value = $g.D
    $ws1.Range("F$r").Value = $g.F
    $ws1.Range("G$r").Value = $g.G
    $ws1.Range("I$r").Value = $g.I
    $ws1.Range("J$r").Value = $g.J
    $ws1.Range("K$r").Value = $g.K
    $ws1.Range("L$r").Value = $g.L

    $mCell = $ws1.Range("M$r")
    $mCell.Formula = "=ABS(D$r-F$r)"
    $mCell.NumberFormat = "#,##0"

    if ($g.Hi) {
        $ws1.Range("K$r").Interior.Color = 5287936
    }
}

# --- Sheet1 sheet view: scrolled down, C142 selected ------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 111
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet2 sheet view: selection moved to R22 (keep Sheet1 as active tab) -
$ws2.Range("R22").Select()
$ws1.Range("C142").Select()

# --- Charts: extend cached source ranges from row 19 to row 20 -------------
foreach ($chartName in @("Chart 1", "Chart 2", "Chart 3")) {
    $co = $ws2.ChartObjects($chartName)
    $chart = $co.Chart
    $series = $chart.SeriesCollection(1)
    $series.XValues = $ws2.Range("A2:A20")
    if ($chartName -eq "Chart 1") {
        $series.Values = $ws2.Range("D2:D20")
    } elseif ($chartName -eq "Chart 2") {
        $series.Values = $ws2.Range("F2:F20")
    } else {
        $series.Values = $ws2.Range("G2:G20")
    }
}

Write-Output "done"
